$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# New location entries (Major City must be interned before Port Village/City
# so the shared-string table order matches: 46 = Major City, 47 = Port Village/City)
$ws.Range("D7").Value = "Major City"
$ws.Range("G10").Value = "Major City"
$ws.Range("I3").Value = "Port Village/City"

# Columns A:I now share the same custom width (was only A:E before).
# NOTE: target OOXML width is 18.7109375; the host's ColumnWidth setter only
# resolves to 1/6-character increments, so 17.8 is the closest input that
# lands on the nearest achievable stored width (18.666666666666668).
$ws.Range("A:I").ColumnWidth = 17.8

# Update the active selection
$ws.Range("I4").Select() | Out-Null
